$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column L: "Statistical Thinking in Python (Part 1)" course ratings ---

# Header cell L1 - same look as the other header cells (J1/K1)
$ws.Range("L1").Value = "Statistical Thinking in Python (Part 1)"
$ws.Range("K1").Copy()
$ws.Range("L1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Data cell L2 - same look as the other data cell K2
$ws.Range("L2").Value = 4
$ws.Range("K2").Copy()
$ws.Range("L2").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Match column width to the sibling style-2 columns (J/K), which are 34.140625
# (OOXML) units wide -- that corresponds to a COM ColumnWidth of raw-5/6.
$ws.Range("L1").EntireColumn.ColumnWidth = 33.307291666666664

# Update the view: the new column is now the active selection, scrolled into view.
$ws.Range("L3").Select() | Out-Null
